$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unlist the existing table (keeps cell values/format, drops table wrapper) ---
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# --- Remove the title row and the blank row above the header ---
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# --- Insert a new column before the "Source" column (col C) for the new "Source" data,
#     pushing the old Source column (which becomes "Link") to column D ---
$ws.Columns.Item(3).Insert()

# --- Header row ---
$ws.Range("C1").Value = "Source"
$ws.Range("D1").Value = "Link"

# --- Row 2: Population dataset ---
$ws.Range("A2").Value = "DCD-area-proypoblacion-Mun-1985-2024-(Population Evolution of Cities)"
$ws.Range("C2").Value = "DANE"
$ws.Range("D2").Value = "https://www.dane.gov.co/index.php/estadisticas-por-tema/demografia-y-poblacion/proyecciones-de-poblacion"
$ws.Range("E2").Value = "1984-2035"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "No"

# --- Row 4: GEIH dataset ---
$ws.Range("A4").Value = "GEIH_cleaned_2015-01-2023-07"
$ws.Range("C4").Value = "DANE"
$ws.Range("D4").Value = "https://www.dane.gov.co/index.php/estadisticas-por-tema/mercado-laboral/empleo-y-desempleo/mercado-laboral-historicos"
$ws.Range("E4").Value = "2015.01-2023.07"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = "Yes (missing previous years)"

# --- Hyperlinks on the Link column ---
$ws.Hyperlinks.Add($ws.Range("D2"), $ws.Range("D2").Value) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), $ws.Range("D4").Value) | Out-Null

# --- Rebuild the table over the new range ---
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:G15"), [System.Type]::Missing, 1)
$lo2.Name = "Tabla1"
$lo2.TableStyle = "TableStyleLight21"

$ws.Range("G4").Select()
